$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Neo4j / DB query text (column B) and the stat query text (column C)
# on every data row (2-4): the recurrence-score filter changes from
# "16-20" to "31-35".
foreach ($r in 2..4) {
    $bCell = $ws.Cells.Item($r, 2)
    $bText = $bCell.Text
    $bCell.Value = $bText.Replace('"16-20"', '"31-35"')

    $cCell = $ws.Cells.Item($r, 3)
    $cText = $cCell.Text
    $cCell.Value = $cText.Replace('"16-20"', '"31-35"')
}

# Update view/selection state: scroll so row 3 is the top-left visible row
# and select C3 (previously D4 was selected with no scroll offset).
$ws.Range("C3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
